$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge "di" + <bookmark _GoBack/> + "fferences." into a single
#    "differences." run, removing the old _GoBack bookmark in the process.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("The main outcome variables of interest are parameters")
$para = $rng.Paragraphs(1).Range
$paraText = $para.Text
$offset = $paraText.IndexOf("differences.")
$diStart = $para.Start + $offset
$bmPos = $diStart + 2

# Minimal edit that crosses the bookmark position so the engine drops it.
$cross = $d.Range($bmPos - 1, $bmPos + 1)
$cross.Text = "IX"
$cross2 = $d.Range($bmPos - 1, $bmPos + 1)
$cross2.Text = "if"

# ---------------------------------------------------------------------------
# 2) Center the figure paragraph and its caption paragraph.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Figure 1: A sample overview over the first few lines of the data") | Out-Null
$captionPara = $rng2.Paragraphs(1)
$captionStart = $captionPara.Range.Start

# the picture paragraph directly precedes the caption paragraph
$precedingRange = $d.Range($captionStart - 1, $captionStart - 1)
$picturePara = $precedingRange.Paragraphs(1)
$picturePara.Format.Alignment = 1
$captionPara.Format.Alignment = 1

# ---------------------------------------------------------------------------
# 3) Move the _GoBack bookmark so that it wraps the caption paragraph:
#    bookmarkStart right before its run content, bookmarkEnd right after the
#    paragraph (i.e. at the start of the following paragraph).
# ---------------------------------------------------------------------------
$bmXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$startIns = $d.Range($captionPara.Range.Start, $captionPara.Range.Start)
$startIns.InsertXML($bmXml)

$rng3 = $d.Content
$rng3.Find.Execute("Figure 1: A sample overview over the first few lines of the data") | Out-Null
$captionPara2 = $rng3.Paragraphs(1)
$afterEnd = $captionPara2.Range.End

$endXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endIns = $d.Range($afterEnd, $afterEnd)
$endIns.InsertXML($endXml)
